$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for every data row (2..101).
# The diff shows every one of these values incrementing from 45174 to 45175
# (i.e. one day later). Update them all by setting the new serial value.
for ($row = 2; $row -le 101; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
